# Adding @Purchase @Tool and updating Database Spreadsheet
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Add the "Purchase" worksheet after the last existing sheet ("School")
# ---------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$purchase = $wb.Worksheets.Add($null, $lastSheet)
$purchase.Name = "Purchase"

$purchase.Range("A1").Value = "Sprocket"
$purchase.Range("B1").Value = 5
$purchase.Range("C1").Value = 40

$purchase.Range("A2").Value = "VEX Motor"
$purchase.Range("B2").Value = 20
$purchase.Range("C2").Value = 30

$purchase.Range("A3").Value = "VEX Brain"
$purchase.Range("B3").Value = 25
$purchase.Range("C3").Value = 10

$purchase.Range("A4").Value = "Angle Gusset"
$purchase.Range("B4").Value = 1
$purchase.Range("C4").Value = 100

# Empty, left-aligned cell below the table
$purchase.Range("C5").HorizontalAlignment = -4131
$purchase.Range("C5").Select()

$purchase.Columns.Item(1).ColumnWidth = 19.5

# ---------------------------------------------------------------------
# Add the "Tools" worksheet after "Purchase"
# ---------------------------------------------------------------------
$tools = $wb.Worksheets.Add($null, $purchase)
$tools.Name = "Tools"

$tools.Range("A1").Value = "Clamp"
$tools.Range("B1").Value = 8

$tools.Range("A2").Value = "Saw"
$tools.Range("B2").Value = 3

$tools.Range("A3").Value = "Phillips Screwdriver"
$tools.Range("B3").Value = 4

$tools.Range("A4").Value = "Punch"
$tools.Range("B4").Value = 5

$tools.Columns.Item(1).ColumnWidth = 15.0

# Match the saved view state: scrolled so row 4 is at the top, selection on D24
$tools.Application.Goto($tools.Range("A4"), $false)
$tools.Range("D24").Select()

# ---------------------------------------------------------------------
# Make "Purchase" the active sheet/tab (matches activeTab=4 in the diff)
# ---------------------------------------------------------------------
$purchase.Activate()
$purchase.Range("C5").Select()
